# Altered to get mean of columns C and D from 3D Cell bodies tab
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column widths: columns 1..31 (A..AE) all become width 20 ---
# This engine stores ColumnWidth + 5/6 as the OOXML "width" attribute,
# so set ColumnWidth to (20 - 5/6) to land on a stored width of exactly 20.
$targetWidth = 20 - (5/6)
for ($col = 1; $col -le 31; $col++) {
    $ws.Columns.Item($col).ColumnWidth = $targetWidth
}

# --- Number format for the "mean" columns E and F (5 and 6) ---
$ws.Range("E1:F2").NumberFormat = "#,##0.00"

# --- Remove old header label in A1 (no replacement header there anymore) ---
$ws.Range("A1").ClearContents()

# --- Row 1 headers (B1:AE1) ---
$ws.Range("B1").Value = "Perimeter(µm)"
$ws.Range("C1").Value = "Area(µm²)"
$ws.Range("D1").Value = "Feret Max(µm)"
$ws.Range("E1").Value = "Enclosed Volume(µm³)"
$ws.Range("F1").Value = "Surface Area(µm²)"
$ws.Range("G1").Value = "n of 3D"
$ws.Range("H1").Value = "Cell Body Length(µm)"
$ws.Range("I1").Value = "Cell Body Mean Length"
$ws.Range("J1").Value = "Cell Body Area (µm²)"
$ws.Range("K1").Value = "Cell Body Mean Area"
$ws.Range("L1").Value = "Cell Body Surface(µm²)"
$ws.Range("M1").Value = "Cell Body Mean Surface"
$ws.Range("N1").Value = "Cell Body Volume(µm³)"
$ws.Range("O1").Value = "Cell Body Mean Volume"
$ws.Range("P1").Value = "Axon Length(µm)"
$ws.Range("Q1").Value = "Axon Mean Length"
$ws.Range("R1").Value = "Axon Area (µm²)"
$ws.Range("S1").Value = "Axon Mean Area"
$ws.Range("T1").Value = "Axon Surface(µm²)"
$ws.Range("U1").Value = "Axon Mean Surface"
$ws.Range("V1").Value = "Axon Volume(µm³)"
$ws.Range("W1").Value = "Axon Mean Volume"
$ws.Range("X1").Value = "Dendrite Length(µm)"
$ws.Range("Y1").Value = "Dendrite Mean Length"
$ws.Range("Z1").Value = "Dendrite Area (µm²)"
$ws.Range("AA1").Value = "Dendrite Mean Area"
$ws.Range("AB1").Value = "Dendrite Surface(µm²)"
$ws.Range("AC1").Value = "Dendrite Mean Surface"
$ws.Range("AD1").Value = "Dendrite Volume(µm³)"
$ws.Range("AE1").Value = "Dendrite Mean Volume"

# --- Row 2 data (A2:AE2) ---
$ws.Range("A2").Value = "05042016 in1 08232019 gfp reconstruction"
$ws.Range("B2").Value = 130.19
$ws.Range("C2").Value = 465.7
$ws.Range("D2").Value = 35.25
$ws.Range("E2").Value = 1063.275887096774
$ws.Range("F2").Value = 524.4065322580645
$ws.Range("G2").Value = 124
$ws.Range("H2").Value = 36868.3
$ws.Range("I2").Value = 38.97
$ws.Range("J2").Value = 96640.7
$ws.Range("K2").Value = 102.16
$ws.Range("L2").Value = 69178.3
$ws.Range("M2").Value = "N/A"
$ws.Range("N2").Value = 3276.01
$ws.Range("O2").Value = "N/A"
$ws.Range("P2").Value = 0
$ws.Range("Q2").Value = 0
$ws.Range("R2").Value = "N/A"
$ws.Range("S2").Value = "N/A"
$ws.Range("T2").Value = 0
$ws.Range("U2").Value = 0
$ws.Range("V2").Value = 0
$ws.Range("W2").Value = 0
$ws.Range("X2").Value = 0
$ws.Range("Y2").Value = 0
$ws.Range("Z2").Value = "N/A"
$ws.Range("AA2").Value = "N/A"
$ws.Range("AB2").Value = 0
$ws.Range("AC2").Value = 0
$ws.Range("AD2").Value = 0
$ws.Range("AE2").Value = 0

Write-Host "Edit complete"
